# Rename the existing sheet and add a new "Confusion Matrix" sheet with
# the SVC classifier's confusion-matrix data.

$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> "Classification Report" -----------------------------
$wsReport = $wb.Worksheets.Item(1)
$wsReport.Name = "Classification Report"

# --- Add the new "Confusion Matrix" sheet right after it -------------------
$wsMatrix = $wb.Worksheets.Add($null, $wsReport)
$wsMatrix.Name = "Confusion Matrix"

# --- Header row (bold, centered, bordered) ---------------------------------
$wsMatrix.Range("B1").Value = "Predicted 1"
$wsMatrix.Range("C1").Value = "Predicted 2"
$wsMatrix.Range("D1").Value = "Predicted 3"

# --- Row labels (bold, centered, bordered) ----------------------------------
$wsMatrix.Range("A2").Value = "Actual 1"
$wsMatrix.Range("A3").Value = "Actual 2"
$wsMatrix.Range("A4").Value = "Actual 3"

# --- Apply the header/label styling (bold font, thin border, center/top) ---
$headerRange = $wsMatrix.Range("B1:D1")
$labelRange = $wsMatrix.Range("A2:A4")

$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$labelRange.Font.Bold = $true
$labelRange.HorizontalAlignment = -4108
$labelRange.VerticalAlignment = -4160
$labelRange.Borders.LineStyle = 1

# --- Confusion-matrix values -------------------------------------------------
$wsMatrix.Range("B2").Value = 6521
$wsMatrix.Range("C2").Value = 4
$wsMatrix.Range("D2").Value = 0

$wsMatrix.Range("B3").Value = 109
$wsMatrix.Range("C3").Value = 7436
$wsMatrix.Range("D3").Value = 38

$wsMatrix.Range("B4").Value = 0
$wsMatrix.Range("C4").Value = 77
$wsMatrix.Range("D4").Value = 151

# --- Restore the originally active sheet (Classification Report) -----------
$wsReport.Activate()
